$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update E9: "Created test spec and merchant model class" -> "Created merchant test spec and merchant model class"
$ws.Range("E9").Value = "Created merchant test spec and merchant model class"

# Update E10 with a new distinct string (was sharing the same string as E9 before the edit)
$ws.Range("E10").Value = "Created tag test spec and tag model class"

# Update the active selection to E9 (matches <selection activeCell="E9" .../>)
$ws.Range("E9").Select()
